# Auto-update Data Telemetría - Ejecución Diaria
# Appends the newest daily snapshot (fecha = 45996) to the "historico" sheet
# and refreshes the "ultimo_snapshot" sheet with that same latest data.

$wb = $excel.ActiveWorkbook

$historico = $wb.Worksheets.Item("historico")
$snapshot  = $wb.Worksheets.Item("ultimo_snapshot")

# The three new rows of data for fecha = 45996 (2025 date serial),
# in column order: fecha, resumen, total_vin, cnt_Conectado 0-2,
# cnt_Intermitente 3-14, cnt_Limitado 15-30+, cnt_Desconectado 31+,
# cnt_Nunca, pct_Conectado 0-2, pct_Intermitente 3-14, pct_Limitado 15-30+,
# pct_Desconectado 31+, pct_Nunca
$newRows = @(
    @(45996, "Telemetría",                      5905, 3557, 493, 188, 663, 1004, 60.24,             8.35, 3.18, 11.23, 17),
    @(45996, "GPS (según REGLA)",                5302, 4677, 341, 103, 174, 7,    88.20999999999999, 6.43, 1.94, 3.28,  0.13),
    @(45996, "GPS (todas con gps_timestamp)",   11200, 9530, 803, 301, 566, 0,    85.09,             7.17, 2.69, 5.05,  0)
)

# --- Append the new rows to the end of "historico" ---
$lastRow = $historico.UsedRange.Rows.Count
foreach ($rowData in $newRows) {
    $lastRow = $lastRow + 1
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $cell = $historico.Cells.Item($lastRow, $col)
        $cell.Value = $rowData[$col - 1]
    }
    # Match the date-formatted style used by the existing "fecha" column cells
    $historico.Cells.Item($lastRow, 1).NumberFormat = "YYYY-MM-DD"
}

# --- Refresh "ultimo_snapshot" with the same latest data (rows 2-4) ---
$row = 1
foreach ($rowData in $newRows) {
    $row = $row + 1
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $cell = $snapshot.Cells.Item($row, $col)
        $cell.Value = $rowData[$col - 1]
    }
    $snapshot.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
}
